# Notification of proposed works to trees in a conservation area - spec sheet update.
# The "File size" row (previously row 16, nested under Documents[] > File) is removed
# from the table, which shifts every subsequent row up by one (rows 17-74 -> 16-73),
# shrinking the used range from A1:I74 to A1:I73 and shifting every merged cell
# range below it up by one row as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Specification")

$ws.Rows.Item(16).Delete()

# Re-assert the two single-cell "merge" ranges (A41 / B41, formerly A42 / B42
# before the shift) that sit on the "Checklist" row - entire-row delete does not
# reliably relocate degenerate 1x1 merges in this runtime.
$ws.Range("A41").Merge()
$ws.Range("B41").Merge()
